$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 620  # H53: was 649.36365
$ws.Cells.Item(53, 9).Value = 524.1667  # I53: was 569.6
$ws.Cells.Item(53, 11).Value = 524.1667  # K53: was 569.6
$ws.Cells.Item(53, 13).Value = 112.8333  # M53: was 67.39999999999998
$ws.Cells.Item(134, 8).Value = 62500  # H134: was 55500
$ws.Cells.Item(134, 10).Value = 62500  # J134: was 55500
$ws.Cells.Item(134, 12).Value = 62500  # L134: was 55500
$ws.Cells.Item(134, 14).Value = -72640  # N134: was -65640
$ws.Cells.Item(135, 8).Value = 985.34375  # H135: was 1032.7273
$ws.Cells.Item(135, 9).Value = 1013.0645  # I135: was 1061.0625
$ws.Cells.Item(135, 11).Value = 9117.5805  # K135: was 9549.5625
$ws.Cells.Item(135, 13).Value = -6582.5805  # M135: was -7014.5625
$ws.Cells.Item(137, 8).Value = 3500.3809  # H137: was 3196.2766
$ws.Cells.Item(137, 9).Value = 2327.92  # I137: was 2107.4285
$ws.Cells.Item(137, 10).Value = 5224.5884  # J137: was 4800.8945
$ws.Cells.Item(137, 11).Value = 6983.76  # K137: was 6322.2855
$ws.Cells.Item(137, 12).Value = 15673.7652  # L137: was 14402.6835
$ws.Cells.Item(137, 13).Value = -4433.76  # M137: was -3772.2855
$ws.Cells.Item(137, 14).Value = -20773.7652  # N137: was -19502.6835
$ws.Cells.Item(138, 8).Value = 2935.7295  # H138: was 2960.8953
$ws.Cells.Item(138, 9).Value = 1269.1052  # I138: was 1317.4445
$ws.Cells.Item(138, 10).Value = 3415.5151  # J138: was 3395.9265
$ws.Cells.Item(138, 11).Value = 3807.3156  # K138: was 3952.3335
$ws.Cells.Item(138, 12).Value = 10246.5453  # L138: was 10187.7795
$ws.Cells.Item(138, 13).Value = 1332.6844  # M138: was 1187.6665
$ws.Cells.Item(138, 14).Value = -20526.5453  # N138: was -20467.7795

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 5214.2  # H63: was 5414.2666
$ws.Cells.Item(63, 9).Value = 4692.75  # I63: was 4801.273
$ws.Cells.Item(63, 10).Value = 7300  # J63: was 7100
$ws.Cells.Item(63, 11).Value = 4692.75  # K63: was 4801.273
$ws.Cells.Item(63, 12).Value = 7300  # L63: was 7100
$ws.Cells.Item(63, 13).Value = -4006.75  # M63: was -4115.273
$ws.Cells.Item(63, 14).Value = -8672  # N63: was -8472
$ws.Cells.Item(66, 8).Value = 5214.2  # H66: was 5414.2666
$ws.Cells.Item(66, 9).Value = 4692.75  # I66: was 4801.273
$ws.Cells.Item(66, 10).Value = 7300  # J66: was 7100
$ws.Cells.Item(66, 11).Value = 23463.75  # K66: was 24006.365
$ws.Cells.Item(66, 12).Value = 36500  # L66: was 35500
$ws.Cells.Item(66, 13).Value = -20031.75  # M66: was -20574.365
$ws.Cells.Item(66, 14).Value = -43364  # N66: was -42364
$ws.Cells.Item(74, 8).Value = 8392987  # H74: was 8971731
$ws.Cells.Item(74, 9).Value = 13892250  # I74: was 15628631
$ws.Cells.Item(74, 11).Value = 13892250  # K74: was 15628631
$ws.Cells.Item(74, 13).Value = -13891376  # M74: was -15627757
$ws.Cells.Item(76, 8).Value = 24998  # H76: was 24999
$ws.Cells.Item(76, 10).Value = 24998  # J76: was 24999
$ws.Cells.Item(76, 12).Value = 24998  # L76: was 24999
$ws.Cells.Item(76, 14).Value = -25674  # N76: was -25675
$ws.Cells.Item(77, 8).Value = 8392987  # H77: was 8971731
$ws.Cells.Item(77, 9).Value = 13892250  # I77: was 15628631
$ws.Cells.Item(77, 11).Value = 69461250  # K77: was 78143155
$ws.Cells.Item(77, 13).Value = -69456882  # M77: was -78138787
$ws.Cells.Item(79, 8).Value = 24998  # H79: was 24999
$ws.Cells.Item(79, 10).Value = 24998  # J79: was 24999
$ws.Cells.Item(79, 12).Value = 24998  # L79: was 24999
$ws.Cells.Item(79, 14).Value = -27338  # N79: was -27339
$ws.Cells.Item(97, 8).Value = 1650.25  # H97: was 1716.3914
$ws.Cells.Item(97, 9).Value = 1757.3182  # I97: was 1906.65
$ws.Cells.Item(97, 10).Value = 472.5  # J97: was 448
$ws.Cells.Item(97, 11).Value = 1757.3182  # K97: was 1906.65
$ws.Cells.Item(97, 12).Value = 472.5  # L97: was 448
$ws.Cells.Item(97, 13).Value = -1261.3182  # M97: was -1410.65
$ws.Cells.Item(97, 14).Value = -1464.5  # N97: was -1440

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3162  # H86: was 3039.8
$ws.Cells.Item(86, 9).Value = 2389.75  # I86: was 2360.9
$ws.Cells.Item(86, 11).Value = 2389.75  # K86: was 2360.9
$ws.Cells.Item(86, 13).Value = -1266.75  # M86: was -1237.9
$ws.Cells.Item(89, 8).Value = 3162  # H89: was 3039.8
$ws.Cells.Item(89, 9).Value = 2389.75  # I89: was 2360.9
$ws.Cells.Item(89, 11).Value = 11948.75  # K89: was 11804.5
$ws.Cells.Item(89, 13).Value = -6332.75  # M89: was -6188.5
$ws.Cells.Item(133, 8).Value = 59233.332  # H133: was 57812.5
$ws.Cells.Item(133, 9).Value = 0  # I133: was 25000
$ws.Cells.Item(133, 10).Value = 59233.332  # J133: was 58870.97
$ws.Cells.Item(133, 11).Value = 0  # K133: was 25000
$ws.Cells.Item(133, 12).ClearContents()  # L133: was 58870.97
$ws.Cells.Item(133, 13).Value = 59233.332  # M133: was -19940
$ws.Cells.Item(133, 14).Value = -69353.33199999999  # N133: was -68990.97

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 0  # H4: was 10000
$ws.Cells.Item(4, 9).Value = 0  # I4: was 10000
$ws.Cells.Item(4, 11).Value = 0  # K4: was 10000
$ws.Cells.Item(4, 13).ClearContents()  # M4: was -9888
$ws.Cells.Item(31, 8).Value = 696576.8  # H31: was 676751.75
$ws.Cells.Item(31, 9).Value = 15683  # I31: was 14871.5625
$ws.Cells.Item(31, 11).Value = 15683  # K31: was 14871.5625
$ws.Cells.Item(31, 13).Value = -15388  # M31: was -14576.5625
$ws.Cells.Item(34, 8).Value = 696576.8  # H34: was 676751.75
$ws.Cells.Item(34, 9).Value = 15683  # I34: was 14871.5625
$ws.Cells.Item(34, 11).Value = 15683  # K34: was 14871.5625
$ws.Cells.Item(34, 13).Value = -15481  # M34: was -14669.5625
$ws.Cells.Item(132, 8).Value = 3670.6155  # H132: was 3849.4167
$ws.Cells.Item(132, 9).Value = 3747.1365  # I132: was 3923.1904
$ws.Cells.Item(132, 10).Value = 3249.75  # J132: was 3333
$ws.Cells.Item(132, 11).Value = 11241.4095  # K132: was 11769.5712
$ws.Cells.Item(132, 12).Value = 9749.25  # L132: was 9999
$ws.Cells.Item(132, 13).Value = -8711.4095  # M132: was -9239.5712
$ws.Cells.Item(132, 14).Value = -14809.25  # N132: was -15059
$ws.Cells.Item(141, 8).Value = 159959.8  # H141: was 181225
$ws.Cells.Item(141, 10).Value = 159959.8  # J141: was 181225
$ws.Cells.Item(141, 12).Value = 159959.8  # L141: was 181225
$ws.Cells.Item(141, 14).Value = -170319.8  # N141: was -191585

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(74, 8).Value = 15000  # H74: was 14200
$ws.Cells.Item(74, 10).Value = 15000  # J74: was 14200
$ws.Cells.Item(74, 12).Value = 45000  # L74: was 42600
$ws.Cells.Item(74, 14).Value = -47122  # N74: was -44722
$ws.Cells.Item(77, 8).Value = 15000  # H77: was 14200
$ws.Cells.Item(77, 10).Value = 15000  # J77: was 14200
$ws.Cells.Item(77, 12).Value = 135000  # L77: was 127800
$ws.Cells.Item(77, 14).Value = -145608  # N77: was -138408
$ws.Cells.Item(87, 8).Value = 2149.25  # H87: was 2598.3333
$ws.Cells.Item(87, 9).Value = 2149.25  # I87: was 2598.3333
$ws.Cells.Item(87, 11).Value = 6447.75  # K87: was 7794.999899999999
$ws.Cells.Item(87, 13).Value = -5199.75  # M87: was -6546.999899999999
$ws.Cells.Item(90, 8).Value = 2149.25  # H90: was 2598.3333
$ws.Cells.Item(90, 9).Value = 2149.25  # I90: was 2598.3333
$ws.Cells.Item(90, 11).Value = 19343.25  # K90: was 23384.9997
$ws.Cells.Item(90, 13).Value = -13103.25  # M90: was -17144.9997
$ws.Cells.Item(140, 8).Value = 233735.08  # H140: was 169372.56
$ws.Cells.Item(140, 9).Value = 233735.08  # I140: was 169372.56
$ws.Cells.Item(140, 11).Value = 701205.24  # K140: was 508117.68
$ws.Cells.Item(140, 13).Value = -696025.24  # M140: was -502937.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 68326.5  # H5: was 78991.8
$ws.Cells.Item(5, 9).Value = 82487.25  # I5: was 104983
$ws.Cells.Item(5, 11).Value = 82487.25  # K5: was 104983
$ws.Cells.Item(5, 13).Value = -82375.25  # M5: was -104871
$ws.Cells.Item(122, 8).Value = 0  # H122: was 287.5
$ws.Cells.Item(122, 9).Value = 0  # I122: was 283.33334
$ws.Cells.Item(122, 10).Value = 0  # J122: was 300
$ws.Cells.Item(122, 11).Value = 0  # K122: was 850.0000200000001
$ws.Cells.Item(122, 12).ClearContents()  # L122: was 900
$ws.Cells.Item(122, 13).ClearContents()  # M122: was 1599.99998
$ws.Cells.Item(122, 14).Value = 0  # N122: was -5800
$ws.Cells.Item(126, 8).Value = 3907.2222  # H126: was 3528.9473
$ws.Cells.Item(126, 9).Value = 3511.5386  # I126: was 3148
$ws.Cells.Item(126, 10).Value = 4936  # J126: was 4957.5
$ws.Cells.Item(126, 11).Value = 10534.6158  # K126: was 9444
$ws.Cells.Item(126, 12).Value = 14808  # L126: was 14872.5
$ws.Cells.Item(126, 13).Value = -8064.6158  # M126: was -6974
$ws.Cells.Item(126, 14).Value = -19748  # N126: was -19812.5
$ws.Cells.Item(132, 8).Value = 76934130  # H132: was 100013784
$ws.Cells.Item(132, 9).Value = 111113310  # I132: was 166668980
$ws.Cells.Item(132, 11).Value = 333339930  # K132: was 500006940
$ws.Cells.Item(132, 13).Value = -333337400  # M132: was -500004410

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 0  # H22: was 662.25
$ws.Cells.Item(22, 10).Value = 0  # J22: was 662.25
$ws.Cells.Item(22, 12).ClearContents()  # L22: was 662.25
$ws.Cells.Item(22, 14).Value = 0  # N22: was -1252.25
$ws.Cells.Item(27, 8).Value = 0  # H27: was 662.25
$ws.Cells.Item(27, 10).Value = 0  # J27: was 662.25
$ws.Cells.Item(27, 12).ClearContents()  # L27: was 662.25
$ws.Cells.Item(27, 14).Value = 0  # N27: was -876.25
$ws.Cells.Item(40, 8).Value = 3364.6086  # H40: was 2995.6072
$ws.Cells.Item(40, 9).Value = 2610.9412  # I40: was 2312.5908
$ws.Cells.Item(40, 11).Value = 2610.9412  # K40: was 2312.5908
$ws.Cells.Item(40, 13).Value = -2474.9412  # M40: was -2176.5908
$ws.Cells.Item(61, 8).Value = 1068.5  # H61: was 1068.6111
$ws.Cells.Item(61, 9).Value = 1068.5  # I61: was 1119.4117
$ws.Cells.Item(61, 10).Value = 0  # J61: was 205
$ws.Cells.Item(61, 11).Value = 1068.5  # K61: was 1119.4117
$ws.Cells.Item(61, 12).Value = 0  # L61: was 205
$ws.Cells.Item(61, 13).ClearContents()  # M61: was -917.4117000000001
$ws.Cells.Item(61, 14).Value = -866.5  # N61: was -609
$ws.Cells.Item(113, 8).Value = 1068.5  # H113: was 1068.6111
$ws.Cells.Item(113, 9).Value = 1068.5  # I113: was 1119.4117
$ws.Cells.Item(113, 10).Value = 0  # J113: was 205
$ws.Cells.Item(113, 11).Value = 1068.5  # K113: was 1119.4117
$ws.Cells.Item(113, 12).Value = 0  # L113: was 205
$ws.Cells.Item(113, 13).ClearContents()  # M113: was 1050.5883
$ws.Cells.Item(113, 14).Value = 1101.5  # N113: was -4545
$ws.Cells.Item(122, 8).Value = 5865.6562  # H122: was 5237.0835
$ws.Cells.Item(122, 9).Value = 5453.923  # I122: was 4762.4194
$ws.Cells.Item(122, 10).Value = 7649.8335  # J122: was 8180
$ws.Cells.Item(122, 11).Value = 16361.769  # K122: was 14287.2582
$ws.Cells.Item(122, 12).Value = 22949.5005  # L122: was 24540
$ws.Cells.Item(122, 13).Value = -13911.769  # M122: was -11837.2582
$ws.Cells.Item(122, 14).Value = -27849.5005  # N122: was -29440
$ws.Cells.Item(133, 8).Value = 0  # H133: was 75000
$ws.Cells.Item(133, 10).Value = 0  # J133: was 75000
$ws.Cells.Item(133, 12).ClearContents()  # L133: was 75000
$ws.Cells.Item(133, 14).Value = 0  # N133: was -80060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 939.8889  # H100: was 857.5454999999999
$ws.Cells.Item(100, 9).Value = 1158.3334  # I100: was 1310.2
$ws.Cells.Item(100, 10).Value = 503  # J100: was 480.33334
$ws.Cells.Item(100, 11).Value = 2316.6668  # K100: was 2620.4
$ws.Cells.Item(100, 12).Value = 1006  # L100: was 960.66668
$ws.Cells.Item(100, 13).Value = -1775.6668  # M100: was -2079.4
$ws.Cells.Item(100, 14).Value = -2088  # N100: was -2042.66668
$ws.Cells.Item(114, 8).Value = 91988.336  # H114: was 99996.664
$ws.Cells.Item(114, 10).Value = 91988.336  # J114: was 99996.664
$ws.Cells.Item(114, 12).Value = 91988.336  # L114: was 99996.664
$ws.Cells.Item(114, 14).Value = -100666.336  # N114: was -108674.664
$ws.Cells.Item(122, 8).Value = 1563.8096  # H122: was 1622.5238
$ws.Cells.Item(122, 9).Value = 1614.7059  # I122: was 1687.2354
$ws.Cells.Item(122, 11).Value = 4844.1177  # K122: was 5061.706200000001
$ws.Cells.Item(122, 13).Value = -2394.1177  # M122: was -2611.706200000001
$ws.Cells.Item(126, 8).Value = 1357.0714  # H126: was 1110.2106
$ws.Cells.Item(126, 9).Value = 1357.0714  # I126: was 1110.2106
$ws.Cells.Item(126, 11).Value = 4071.2142  # K126: was 3330.6318
$ws.Cells.Item(126, 13).Value = -1601.2142  # M126: was -860.6318000000001
$ws.Cells.Item(132, 8).Value = 2671.9048  # H132: was 2765.7896
$ws.Cells.Item(132, 9).Value = 2708  # I132: was 2817.1765
$ws.Cells.Item(132, 11).Value = 8124  # K132: was 8451.529500000001
$ws.Cells.Item(132, 13).Value = -5594  # M132: was -5921.529500000001
